# "more work done on project"
#
# Adds a second data column (Comparison-2 RMSE / error-rate) to the
# "Given & Regression Temp RMSE" / "Error (RMSE/ Given Temp Range)" summary
# rows at the bottom of Sheet1, inserts a blank separator row above them,
# nudges the saved cell-selection, and sets the sheet to print in portrait
# orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New blank separator row (row 21), matching the look of the other
#     blank separator rows (4, 9, 15) but with a single uniform style
#     across every column A:H (same style already used by A4/A9/A15/etc.). ---
$ws.Range("A4").Copy()
$ws.Range("A21:H21").PasteSpecial(-4122)   # xlPasteFormats

# --- B22: new "Comparison 2" RMSE figure, formatted like its row-mate A22. ---
$ws.Range("A22").Copy()
$ws.Range("B22").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("B22").Value = 5.2831

# --- B23: new "Comparison 2" error-rate formula (RMSE / given-temp-range),
#     percentage-formatted like its row-mates C23:H23, but left-aligned
#     (format seeded from A2, which already carries the right border/no
#     alignment, then switched to the 0.00% number format). ---
$ws.Range("A2").Copy()
$ws.Range("B23").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("B23").NumberFormat = "0.00%"
$ws.Range("B23").Formula = "=B22/24"

$excel.CutCopyMode = $false

# --- Page setup: print portrait. ---
$ws.PageSetup.Orientation = 1              # xlPortrait

# --- Move the saved selection cursor. ---
$ws.Range("L12").Select() | Out-Null
